# Updates cryptos list: refresh price (D) and volume-change (E) columns,
# and reorder a few coin rows (B/C/D/E) to match the latest ranking snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.260.30'
$ws.Range("E2").Value = '  -2.05%  '
$ws.Range("D3").Value = '3.380.78'
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.85'
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.29'
$ws.Range("E6").Value = '  -6.27%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.380.13'
$ws.Range("E8").Value = '  -1.88%  '
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("E10").Value = '  -4.81%  '
$ws.Range("E11").Value = '  -1.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.387'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").Value = '3.956.99'
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.95'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '3.384.03'
$ws.Range("E16").Value = '  -1.64%  '
$ws.Range("E17").Value = '  -3.39%  '
$ws.Range("D18").Value = '60.408.53'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.20'
$ws.Range("E19").Value = '  -1.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.93'
$ws.Range("E20").Value = '  -2.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.00'
$ws.Range("E21").Value = '  -5.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '385.77'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.556'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.20'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -6.08%  '
$ws.Range("D27").Value = '3.528.64'
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.36'
$ws.Range("E30").Value = '  -5.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -4.33%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("E33").Value = '  -7.71%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.57'
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("D36").Value = '3.411.55'
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.88'
$ws.Range("E37").Value = '  -2.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '167.31'
$ws.Range("E38").Value = '  +0.53%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  -4.96%  '
$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.91'
$ws.Range("E40").Value = '  -7.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0769'
$ws.Range("E41").Value = '  -2.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.86'
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.780'
$ws.Range("E43").Value = '  -1.60%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.21'
$ws.Range("E47").Value = '  -2.53%  '
$ws.Range("D48").Value = '2.513.08'
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.11'
$ws.Range("E49").Value = '  -4.24%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.99'
$ws.Range("E50").Value = '  -1.09%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.77'
$ws.Range("E51").Value = '  -3.82%  '
